$wb = $excel.ActiveWorkbook

# --- Insert the new "February 2023" sheet as the first tab ---
$feb = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$feb.Name = "February 2023"

# --- Populate header row ---
$feb.Range("A1").Value = "Mark as Good (Remove BadCIS)"
$feb.Range("B1").Value = "Mark as Bad (add BadCIS)"

# --- Populate data (column A: newly marked-good CIS; column B: newly marked-bad CIS) ---
$feb.Range("A2").Value = "LBC_m44_c34_highgain"

$badList = @(
    "LBA_m02_c06_lowgain",
    "LBA_m02_c06_highgain",
    "LBA_m62_c26_highgain",
    "LBA_m35_c08_highgain",
    "LBA_m38_c46_lowgain",
    "LBC_m62_c08_highgain",
    "EBA_m16_c17_highgain",
    "EBC_m34_c41_lowgain",
    "LBA_m30_c13_lowgain",
    "LBC_m19_c26_highgain",
    "LBC_m57_c06_highgain",
    "EBA_m40_c35_highgain",
    "EBC_m20_c10_lowgain",
    "EBC_m23_c01_lowgain",
    "EBC_m34_c38_highgain",
    "LBA_m03_c17_lowgain",
    "LBA_m51_c12_highgain"
)

$row = 2
foreach ($item in $badList) {
    $feb.Cells.Item($row, 2).Value = $item
    $row = $row + 1
}

$feb.Range("A3").Value = "EBC_m20_c31_highgain"

Write-Host "Done"
